$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1168.5714
$ws.Range("I33").Value = 599.2
$ws.Range("J33").Value = 2592
$ws.Range("K33").Value = 599.2
$ws.Range("L33").Value = 2592
$ws.Range("M33").Value = -370.2
$ws.Range("N33").Value = -3050

$ws.Range("H88").Value = 1407
$ws.Range("J88").Value = 1420.3636
$ws.Range("L88").Value = 1420.3636
$ws.Range("N88").Value = -2232.3636

$ws.Range("H91").Value = 1407
$ws.Range("J91").Value = 1420.3636
$ws.Range("L91").Value = 1420.3636
$ws.Range("N91").Value = -4228.3636

$ws.Range("H92").Value = 185.8421
$ws.Range("I92").Value = 107.76471
$ws.Range("J92").Value = 849.5
$ws.Range("K92").Value = 107.76471
$ws.Range("L92").Value = 849.5
$ws.Range("M92").Value = 1140.23529
$ws.Range("N92").Value = -3345.5

$ws.Range("H100").Value = 2095.2
$ws.Range("I100").Value = 1359
$ws.Range("K100").Value = 1359
$ws.Range("M100").Value = -818

$ws.Range("H127").Value = 2153.5
$ws.Range("I127").Value = 792.25
$ws.Range("K127").Value = 2376.75
$ws.Range("M127").Value = 2583.25

$ws.Range("H132").Value = 2903.1035
$ws.Range("I132").Value = 2863.9285
$ws.Range("K132").Value = 8591.7855
$ws.Range("M132").Value = -6061.7855

$ws.Range("H137").Value = 3056.3635
$ws.Range("I137").Value = 2592.9666
$ws.Range("K137").Value = 7778.899800000001
$ws.Range("M137").Value = -5228.899800000001

$ws.Range("H138").Value = 4124.099
$ws.Range("J138").Value = 4167.2354
$ws.Range("L138").Value = 12501.7062
$ws.Range("N138").Value = -22781.7062

$ws.Range("H141").Value = 3505.158
$ws.Range("I141").Value = 2441.75
$ws.Range("J141").Value = 4278.5454
$ws.Range("K141").Value = 7325.25
$ws.Range("L141").Value = 12835.6362
$ws.Range("M141").Value = -2145.25
$ws.Range("N141").Value = -23195.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7885.442
$ws.Range("I61").Value = 7560.816
$ws.Range("J61").Value = 10352.6
$ws.Range("K61").Value = 7560.816
$ws.Range("L61").Value = 10352.6
$ws.Range("M61").Value = -7348.816
$ws.Range("N61").Value = -10776.6

$ws.Range("H74").Value = 5764.4
$ws.Range("I74").Value = 3109.318
$ws.Range("K74").Value = 3109.318
$ws.Range("M74").Value = -2235.318

$ws.Range("H77").Value = 5764.4
$ws.Range("I77").Value = 3109.318
$ws.Range("K77").Value = 15546.59
$ws.Range("M77").Value = -11178.59

$ws.Range("H136").Value = 7885.442
$ws.Range("I136").Value = 7560.816
$ws.Range("J136").Value = 10352.6
$ws.Range("K136").Value = 22682.448
$ws.Range("L136").Value = 31057.8
$ws.Range("M136").Value = -20132.448
$ws.Range("N136").Value = -36157.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4304.123
$ws.Range("I134").Value = 3448.0544
$ws.Range("J134").Value = 9012.5
$ws.Range("K134").Value = 10344.1632
$ws.Range("L134").Value = 27037.5
$ws.Range("M134").Value = -7809.163199999999
$ws.Range("N134").Value = -32107.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8872.973
$ws.Range("I31").Value = 4533.905
$ws.Range("K31").Value = 4533.905
$ws.Range("M31").Value = -4238.905

$ws.Range("H34").Value = 8872.973
$ws.Range("I34").Value = 4533.905
$ws.Range("K34").Value = 4533.905
$ws.Range("M34").Value = -4331.905

$ws.Range("H62").Value = 3600
$ws.Range("I62").Value = 3600
$ws.Range("K62").Value = 3600
$ws.Range("M62").Value = -2976

$ws.Range("H65").Value = 3600
$ws.Range("I65").Value = 3600
$ws.Range("K65").Value = 18000
$ws.Range("M65").Value = -14880

$ws.Range("H99").Value = 12453.2
$ws.Range("I99").Value = 12453.2
$ws.Range("K99").Value = 12453.2
$ws.Range("M99").Value = -10955.2

$ws.Range("H126").Value = 12453.2
$ws.Range("I126").Value = 12453.2
$ws.Range("K126").Value = 37359.60000000001
$ws.Range("M126").Value = -34889.60000000001

$ws.Range("H134").Value = 3754.7576
$ws.Range("I134").Value = 3200.1538
$ws.Range("K134").Value = 9600.4614
$ws.Range("M134").Value = -7065.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4857.4375
$ws.Range("I102").Value = 3917.6924
$ws.Range("J102").Value = 8929.666999999999
$ws.Range("K102").Value = 3917.6924
$ws.Range("L102").Value = 8929.666999999999
$ws.Range("M102").Value = -2295.6924
$ws.Range("N102").Value = -12173.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5390.913
$ws.Range("I100").Value = 5094.619
$ws.Range("K100").Value = 5094.619
$ws.Range("M100").Value = -4553.619

$ws.Range("H122").Value = 4473.2812
$ws.Range("I122").Value = 3771.182
$ws.Range("J122").Value = 6017.9
$ws.Range("K122").Value = 11313.546
$ws.Range("L122").Value = 18053.7
$ws.Range("M122").Value = -8863.545999999998
$ws.Range("N122").Value = -22953.7

$ws.Range("H136").Value = 6309.1714
$ws.Range("I136").Value = 2755.2
$ws.Range("J136").Value = 8974.65
$ws.Range("K136").Value = 8265.599999999999
$ws.Range("L136").Value = 26923.95
$ws.Range("M136").Value = -5715.599999999999
$ws.Range("N136").Value = -32023.95

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 88212
$ws.Range("J46").Value = 88212
$ws.Range("L46").Value = 88212
$ws.Range("N46").Value = -88674

$ws.Range("H75").Value = 91428.42999999999
$ws.Range("I75").Value = 93333
$ws.Range("K75").Value = 93333
$ws.Range("M75").Value = -92397

$ws.Range("H78").Value = 91428.42999999999
$ws.Range("I78").Value = 93333
$ws.Range("K78").Value = 279999
$ws.Range("M78").Value = -275319

$ws.Range("H81").Value = 1249.5
$ws.Range("I81").Value = 1266
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 2532
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -1471
$ws.Range("N81").Value = -4522

$ws.Range("H84").Value = 1249.5
$ws.Range("I84").Value = 1266
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 12660
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -7356
$ws.Range("N84").Value = -22608

$ws.Range("H96").Value = 33749.5
$ws.Range("I96").Value = 24999
$ws.Range("K96").Value = 24999
$ws.Range("M96").Value = -23626

$ws.Range("H132").Value = 2259.8572
$ws.Range("J132").Value = 3747.625
$ws.Range("L132").Value = 11242.875
$ws.Range("N132").Value = -16302.875

$ws.Range("H134").Value = 88212
$ws.Range("J134").Value = 88212
$ws.Range("L134").Value = 264636
$ws.Range("N134").Value = -269706
